$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (styles) of the last existing data row (row 15) down
# into the new row 16 so the new row inherits the same look (bold/border
# index column, date-formatted column, etc.) before we populate values.
$ws.Range("A15:G15").Copy()
$ws.Range("A16:G16").PasteSpecial(-4122)  # xlPasteFormats

# New data row for 2020-06-15 (serial 43997)
$ws.Range("A16").Value = 14
$ws.Range("B16").Value = 43997
$ws.Range("C16").Value = 150264
$ws.Range("D16").Value = 211616
$ws.Range("E16").Value = 53217
$ws.Range("F16").Value = 17580
$ws.Range("G16").Value = 32.36

# B15 switches from the date-only format to the date+time format, while
# B16 (the newly inserted row) takes on the date-only format that B15
# used to have.
$ws.Range("B15").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B16").NumberFormat = "YYYY-MM-DD"
